# Apply updates to column C (the "words" column) for specific rows,
# identified by their key in column A, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "FW" = "két, golf, Merdeka, gas, bar, fee, judo, bêtông, tivi, ôtô"
    "I"  = "ạ, ơi, Ôi, Thôi, à, Ối, Trời ơi, nha, nhỉ, hỡi"
    "Np" = "Văn, VN, Mai, Nguyễn, Hoà, Trần, Sài Gòn, Nhật, Tiết, Minh"
    "Ny" = "HCV, HLV, TP., UBND, VĐV, Q., NV, P., TP, CLB"
    "T"  = "cả, thôi, chính, đấy, mà, tận, cái, ngay, Chính, trời"
    "X"  = "như thế, như vậy, nhất là, làm sao, có lẽ, làm gì, Vậy mà, hầu hết, Vì vậy, Thật ra"
    "Z"  = "bất, phi, siêu, gia, đại, liên, bán, phó"
}

$dim = $ws.UsedRange
$rowCount = $dim.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $key = $ws.Cells.Item($r, 1).Value2
    if ($null -ne $key -and $updates.ContainsKey([string]$key)) {
        $ws.Cells.Item($r, 3).Value2 = $updates[[string]$key]
    }
}

$wb.Save()
